$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B->C, old C->D) to hold the new
# "StatQuery" data.
$ws.Range("B1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "StatQuery"

# New stat-query text for row 2 of the inserted column, using the same
# wrap-text style already used by A2 (style index 1 / wrapText).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Salivary gland cancer']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# New column B should be as wide as column A (75.81640625 chars). The
# engine quantizes ColumnWidth to 1/6-character steps, so 75 is the input
# that lands closest (75.8333...) to column A's exact stored width.
$ws.Range("B1").EntireColumn.ColumnWidth = 75

# Selection as shown in the diff.
$ws.Range("A2").Select()
